# Generate Report for Handoff
# - Row for "7a3b1b1b-...md" moves from "Handed back: in sync with en-US" to
#   "Ready for handoff" on the Overview sheet as well as the zh-cn / de-de
#   per-language detail sheets.
# - The corresponding "Latest Handoff Datetime" values for that batch are
#   refreshed on the zh-cn and de-de sheets.

$wb = $excel.ActiveWorkbook

$status = "Ready for handoff"

# --- Overview sheet --------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B3").Value2 = $status
$overview.Range("C3").Value2 = $status

# --- zh-cn detail sheet ------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("B3").Value2 = $status
$zhcn.Range("D2").Value2 = "2016-03-07 02:52:35"
$zhcn.Range("D3").Value2 = "2016-03-07 02:52:35"

# --- de-de detail sheet ------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("B3").Value2 = $status
$dede.Range("D2").Value2 = "2016-03-07 02:52:45"
$dede.Range("D3").Value2 = "2016-03-07 02:52:45"
